# Update public/menu.xlsx data: replace the old menu rows with the new
# weekly menu (07.04.2025), extend the table through row 10, add a
# trailing row 23, and give the Date / Price columns their own number
# formats (date + PLN currency).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "d/mm/yyyy"
$currencyFormat = '#,##0.00\ [$zł-415];[RED]\-#,##0.00\ [$zł-415]'
$dateText = "07.04.2025"

# Wipe out the old menu rows (2-5, cols A-D) completely - content AND
# formatting - so nothing lingers from the previous Image Path column.
$ws.Range("A2:D5").Clear()

# row -> (use date-format style on col A?, dish name, price)
$menu = @(
    @{ Row = 2;  DateStyled = $true;  Name = "Gulasz wieprzowy, kasza gryczana, surówka + zupa wiosenna lub grochowa "; Price = 32 },
    @{ Row = 3;  DateStyled = $false; Name = "Grillowany rumsztyk z cebulką, ziemniaki, surówka + zupa wiosenna lub grochowa "; Price = 31 },
    @{ Row = 4;  DateStyled = $false; Name = "Gulasz wieprzowy, kasza gryczana, surówka"; Price = 29 },
    @{ Row = 5;  DateStyled = $true;  Name = "Grillowany rumsztyk z cebulką, ziemniaki, surówka "; Price = 29 },
    @{ Row = 6;  DateStyled = $false; Name = "Zupa wiosenna lub grochowa "; Price = 9.5 },
    @{ Row = 7;  DateStyled = $true;  Name = "Kotlet schabowy, ziemniaki, surówka"; Price = 28 },
    @{ Row = 8;  DateStyled = $true;  Name = "Panierowany filet z kurczaka, ziemniaki, surówka"; Price = 28 },
    @{ Row = 9;  DateStyled = $true;  Name = "Bigos z ziemniakami "; Price = 24 },
    @{ Row = 10; DateStyled = $false; Name = "Naleśniki z serkiem i brzoskwiniami "; Price = 21 }
)

foreach ($item in $menu) {
    $r = $item.Row

    $aCell = $ws.Range("A$r")
    if ($item.DateStyled) {
        $aCell.NumberFormat = $dateFormat
    }
    # Assign through a text formula, then flatten the formula result back
    # into the cell as a plain value (copy / paste-special-values) - this
    # stores the date as literal text (matching the source data) instead
    # of letting Excel coerce it into a date serial number, while leaving
    # whatever number format is already on the cell untouched.
    $aCell.Formula = '="' + $dateText + '"'
    $aCell.Copy()
    $aCell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $bCell = $ws.Range("B$r")
    $bCell.Value = $item.Name

    $cCell = $ws.Range("C$r")
    $cCell.NumberFormat = $currencyFormat
    $cCell.Value = $item.Price
}

# Trailing row with a single blank-ish cell.
$ws.Range("B23").Value = " "

# Match the author's last selection.
[void]$ws.Range("C11").Select()
